# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as row 50 of the sheet
# (pushing the previously-existing rows 50..121 down to 51..122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 50; Excel shifts rows 50..121 down to
# 51..122 and the new row inherits formatting from the row above (this is
# what gives the new D50 cell the date-style "s=2" seen in the original
# file on every D column cell).
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly observation.
$ws.Cells.Item(50, 1).Value  = 11
$ws.Cells.Item(50, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value  = "Bíobío"
$ws.Cells.Item(50, 4).Value  = 44540
$ws.Cells.Item(50, 5).Value  = 8
$ws.Cells.Item(50, 6).Value  = "Fruta"
$ws.Cells.Item(50, 7).Value  = 100108
$ws.Cells.Item(50, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(50, 9).Value  = 100108005
$ws.Cells.Item(50, 10).Value = "Piña"
$ws.Cells.Item(50, 11).Value = "Caramelo"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 150
$ws.Cells.Item(50, 14).Value = 17000
$ws.Cells.Item(50, 15).Value = 18000
$ws.Cells.Item(50, 16).Value = 17533
$ws.Cells.Item(50, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(50, 18).Value = "Ecuador"
$ws.Cells.Item(50, 19).Value = 1252
$ws.Cells.Item(50, 20).Value = 14
